$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: fill in the new log entry (previously blank)
$ws.Range("A20").Value = "Replaced segue to detail view"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 44621
$ws.Range("D20").Value = "Replaced segue to detailview with new method storyboard.instantiateViewController. Removed optional HouseManager (view model). Made code a bit cleaner."

# Row 30: the "Total amount of hours" label gets pushed out to a new shared
# string (since two new strings were inserted before it). Re-set it so the
# shared string table order is rebuilt to match, and let the SUMIF formula
# recalc pick up the new hours total automatically.
$ws.Range("A30").Value = "Total amount of hours"
